# Add two new "timeout" configuration columns to the ServerDatabaseConfig
# sheet: ConnectionTimeOut (30) and QueryTimeOut (60), inserted right after
# DatabaseName and before the existing WorkersCount column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ServerDatabaseConfig")

# Remember the width of the neighbouring column so the two freshly
# inserted columns pick up a matching (header-row) look & feel.
$neighborWidth = $ws.Range("B1").ColumnWidth

# Insert two blank columns at C:D (WorkersCount, formerly in C, shifts to E).
$null = $ws.Range("C1:D1").EntireColumn.Insert()

# Give the new columns a similar width to column B.
$ws.Range("C1:D1").ColumnWidth = $neighborWidth

# Header row.
$ws.Range("C1").Value = "ConnectionTimeOut"
$ws.Range("D1").Value = "QueryTimeOut"

# Data row.
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 60

# Match the saved selection from the source workbook.
$null = $ws.Range("C6").Select()
